# Apply "Playing with XP numbers" edits to lw2_xpdata.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Current" (2nd tab) ---
$wsCurrent = $wb.Worksheets.Item("Current")

$wsCurrent.Range("G2").Value = 5
$wsCurrent.Range("H2").Value = 5.9
$wsCurrent.Range("H4").Value = 5.8
$wsCurrent.Range("G5").Value = 4
$wsCurrent.Range("G8").Value = 5
$wsCurrent.Range("H9").Value = 5.6

# Update the active cell / selection on the "Current" sheet
$wsCurrent.Activate()
$wsCurrent.Range("G7").Select()

# --- Sheet "Rebalance" (3rd tab) ---
$wsRebalance = $wb.Worksheets.Item("Rebalance")
$wsRebalance.Activate()
$wsRebalance.Range("G2:H9").Select()

# Leave "Current" as the active/selected sheet at the end, matching tabSelected in the diff
$wsCurrent.Activate()
